$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A1").Value = "Biokart India Pvt Ltd – Bringing Biotech Closer!"
